$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K11").Value = -0.3352267436446591
$ws.Range("J12").Value = 0.01855976243503714
$ws.Range("I13").Value = -0.1296176279974082
$ws.Range("H14").Value = -0.2870636170015632
$ws.Range("G15").Value = 0.2135958395245076
$ws.Range("F16").Value = -0.06676204101096155
$ws.Range("E17").Value = 0.1052128168340501
$ws.Range("D18").Value = -0.2006497229122814
$ws.Range("C19").Value = 0.4116802297750048
$ws.Range("B20").Value = -0.2766911554241067
